$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix the filenames in column B for the ToxA rows (remove the dash: "Tox-A" -> "ToxA")
$ws.Range("B4").Value = "ToxA_R1_98_S2_L008_R1_001_x.fastq.gz"
$ws.Range("B5").Value = "ToxA_R2_S2_L005_R1_001_x.fastq.gz"

# Row 5's A value becomes numeric 3 (previously blank)
$ws.Range("A5").Value = 3

# Remove the ToxB rows (old rows 6 and 7) entirely
$ws.Rows("6:7").Delete()
